$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (this also updates the <sheet name="..."> in workbook.xml)
$ws.Name = "Through 2022-06-10"

# Update the row label for June
$ws.Range("A7").Value = "June (through 06-10)"

# Update the June row (row 7) values
$ws.Range("B7").Value = 2
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 11
$ws.Range("G7").Value = 42
$ws.Range("H7").Value = 35
$ws.Range("I7").Value = 35

# Update the Total row (row 8) values
$ws.Range("B8").Value = 110
$ws.Range("E8").Value = 315
$ws.Range("F8").Value = 215
$ws.Range("G8").Value = 400
$ws.Range("H8").Value = 666
$ws.Range("I8").Value = 698
